$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.184.87'
$ws.Range("E2").Value = '  +1.86%  '
$ws.Range("D3").Value = '2.052.93'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.29%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.06'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.382'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.40'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0756'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.44%  '
$ws.Range("E12").Value = '  +1.29%  '
$ws.Range("D13").Value = '2.355.33'
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.91'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.772'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.13%  '
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").Value = '2.050.45'
$ws.Range("D19").Value = '37.134.51'
$ws.Range("E19").Value = '  +1.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +9.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '68.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("D22").Value = '0.0₃0808'
$ws.Range("E22").Value = '  +1.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '225.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.12%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.19%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.38%  '
$ws.Range("E30").Value = '  +0.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.126'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.45'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0615'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.57'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.51'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.05%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.27'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.74'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.71'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.95'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.45'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").Value = '1.482.14'
$ws.Range("E43").Value = '  +0.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.41%  '
$ws.Range("E45").Value = '  +5.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0927'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0210'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.25%  '
$ws.Range("E48").Value = '  +1.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '15.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.17'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.94'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.75%  '
